# fix: conserta erro e comenta o código
# Corrige valores incorretos na planilha de dados do Reclame Aqui.
# Os valores são mantidos como texto (não numéricos), pois assim
# estão armazenados originalmente na planilha.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Garante que as células continuem formatadas/tratadas como texto,
# evitando que o Excel converta os percentuais/números automaticamente.
$cells = @("D2", "C4", "D4", "F5", "D7", "E7", "F7")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Lojas Renner - "Voltariam a fazer negócio": 73.5% -> 73.4%
$ws.Range("D2").Value = "73.4%"

# Netshoes - "Reclamações Respondidas": 97.9% -> 98%
$ws.Range("C4").Value = "98%"
# Netshoes - "Voltariam a fazer negócio": 72.8% -> 72.7%
$ws.Range("D4").Value = "72.7%"

# Yeesco - "Nota do Consumidor": 2.75 -> 2.74
$ws.Range("F5").Value = "2.74"

# Wanted - "Voltariam a fazer negócio": 31.5% -> 31.4%
$ws.Range("D7").Value = "31.4%"
# Wanted - "Índice de Solução": 60.9% -> 60.5%
$ws.Range("E7").Value = "60.5%"
# Wanted - "Nota do Consumidor": 3.28 -> 3.26
$ws.Range("F7").Value = "3.26"
